# The "Squad Total" footer row (row 35) had its summary formulas/values wiped
# out in the committed workbook: every cell in that row is now empty, but the
# cells keep their original number formats/styles (only the cell contents were
# removed, not the row itself and not the formatting).
#
# This matches a user selecting the entire row 35 and pressing Delete, which
# is reproduced here as "select the row, then clear its contents".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = $ws.Rows("35:35")
$targetRow.Select()
$targetRow.ClearContents()
